$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the calibration job data that had been recorded in columns BF:BN
# for rows 34-42 (mean/std no., DUT datetime, mean/std RH, mean/std T,
# mean/std Td/f). The cells become blank, matching the already-blank
# rows 43:44 below them.
$range = $ws.Range("BF34:BN42")
$range.ClearContents()

# The DUT datetime column (BH) carried a date/time number format; clear
# that formatting too so the now-empty cells fall back to General,
# consistent with the other blank rows in this block.
$ws.Range("BH34:BH42").ClearFormats()
